# Added comments to work plan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

$ws.Range("G6").Value = "Stworzenie szkieletu Computational Noda(3h)."
$ws.Range("G5").Value = "Zapoznanie się z dokumentacją, wyjaśnienie niejasność. Podział prac w grupie(3h)."
$ws.Range("G7").Value = "Projektowanie i implementacja Computational Noda(3h)."
$ws.Range("G8").Value = "Projektowanie i implementacja Computational Noda(3h)."
$ws.Range("G9").Value = "Wprowadzenie poprawek zwiazanych ze zmianą sposobu komunikacji(dodatkowe 4h)."
$ws.Range("G10").Value = "Refaktoryzacja kodu."

# Update the view state to match the saved workbook (scrolled up one row,
# selection moved from E10 to G10)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("G10").Select()
